$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reformat the sample-year text in row 3 (SampleYr column, P3): "2002 to 2008" -> "2002-2008"
$ws.Range("P3").Value = "2002-2008"

# Move the active selection from R14 to P4
[void]$ws.Range("P4").Select()
